$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.840.94'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '3.490.74'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  +3.67%  '
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').Value = '4.095.34'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '66.837.43'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '3.499.77'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '394.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000121'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('E30').Value = '  -2.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '162.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('E38').Value = '  +3.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.13'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').Value = '2.821.76'
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('E46').Value = '  -2.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '335.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '34.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.08%  '
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  -1.14%  '
